# Fruta / hortaliza, semanal
#
# Insert a new weekly record as the new row 16 ("Hortaliza, Macroferia
# Regional de Talca - Espárragos"), pushing the previously existing rows
# 16-30 down to 17-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 16; this shifts rows 16:30
# down to 17:31 (values, formats and styles move with them).
$ws.Rows("16:16").Insert()

# Populate the newly inserted row 16 with this week's data.
$ws.Cells.Item(16, 1).Value = 5
$ws.Cells.Item(16, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(16, 3).Value = "Maule"
$ws.Cells.Item(16, 4).Value = 44484
$ws.Cells.Item(16, 5).Value = 7
$ws.Cells.Item(16, 6).Value = 300000000
$ws.Cells.Item(16, 7).Value = "Espárragos"
$ws.Cells.Item(16, 8).Value = "Verde"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 5000
$ws.Cells.Item(16, 11).Value = 800
$ws.Cells.Item(16, 12).Value = 900
$ws.Cells.Item(16, 13).Value = 840
$ws.Cells.Item(16, 14).Value = "`$/kilo"
$ws.Cells.Item(16, 15).Value = "Provincia de Linares"
$ws.Cells.Item(16, 16).Value = 840
$ws.Cells.Item(16, 17).Value = 1
$ws.Cells.Item(16, 18).Value = "Hortaliza"
